$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: row "3" (the row keyed r="3" in the sheet XML) becomes what row "4" used to hold ---
# Copying preserves the original (text) cell typing instead of re-inferring types from scratch.
$srcRow4 = $ws.Range($ws.Cells.Item(4,1), $ws.Cells.Item(4,6))
$dstRow3 = $ws.Range($ws.Cells.Item(3,1), $ws.Cells.Item(3,6))
$srcRow4.Copy($dstRow3)

# --- Step 2: row "4" gets brand-new subscriber data ---
$ws.Cells.Item(4,1).Value2 = "naoures"
$ws.Cells.Item(4,2).Value2 = "bzeouich"
# these two look numeric, force them to stay text (as every other id/phone column in this sheet is)
$ws.Cells.Item(4,3).Value2 = "'14034112"
$ws.Cells.Item(4,4).Value2 = "'97944447"
$ws.Cells.Item(4,5).Value2 = "MOKNINE"
$ws.Cells.Item(4,6).Value2 = "النحل"

# --- Step 3: drop the two trailing rows (old r="5" and r="6") entirely ---
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
